$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "MERSEN "
$ws.Range("B50").Value = "MERSEN"
$ws.Range("C50").Value = "FER"
$ws.Range("E50").Value = "OUI"

$ws.Range("C51").Value = "SCH"
$ws.Range("A51").Value = "SCHNEIDER ELECTRIC"
$ws.Range("B51").Value = "SCHNEIDER ELECTRIC"
$ws.Range("E51").Value = "OUI"

$null = $ws.Range("F40").Select()
